# Form-4-Workplan.docx fix: force a fixed table layout on the document's
# tables (adds <w:tblLayout w:type="fixed"/> to each table's <w:tblPr>,
# right after <w:tblBorders> and before <w:tblLook>, matching Word's own
# placement order for that element).
$d = $word.ActiveDocument

foreach ($t in $d.Tables) {
    # Turning off Word's "AutoFit" behavior is what persists
    # <w:tblLayout w:type="fixed"/> into the table properties - this stops
    # the table from resizing itself to fit its contents/window, which is
    # the actual bug being fixed (columns shifting around).
    $t.AutoFitBehavior(0)
}
